$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.003208871385164791
$ws.Range("C2").Value = 10.34677158129881
$ws.Range("D2").Value = 57107556.33100624
$ws.Range("E2").Value = 85231193291209616
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 85231193348317184
